$wb = $excel.ActiveWorkbook

# Add a new worksheet named "Sheet1" positioned after the existing "NCAP_BND" sheet
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "Sheet1"

# Populate the new sheet with the DINS transform data
$ws.Range("B6").Value = "~TFM_DINS-TS"

$ws.Range("B7").Value = "attribute"
$ws.Range("C7").Value = "process"
$ws.Range("D7").Value = 2000

$ws.Range("B8").Value = "act_cost"
$ws.Range("C8").Value = "ELCNENUC00"
$ws.Range("D8").Value = 0.44

# Auto-fit the used columns to their content, like Excel does on entry
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null

# Make the new sheet the active sheet/tab, with I18 selected (matches target selection)
$ws.Activate()
[void]$ws.Range("I18").Select()
